$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = '@'
    $c.Value = $val
    $c.Style = 'Normal'
}

$ws.Range('D2').Value = '67.607.48'
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').Value = '3.319.49'
$ws.Range('E3').Value = '  +2.10%  '
$ws.Range('E4').Value = '  +0.09%  '
Set-TextValue 'D5' '580.37'
$ws.Range('E5').Value = '  +0.43%  '
Set-TextValue 'D6' '174.07'
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +2.30%  '
$ws.Range('D9').Value = '3.314.98'
$ws.Range('E9').Value = '  +1.97%  '
$ws.Range('E10').Value = '  +5.49%  '
$ws.Range('E11').Value = '  +1.67%  '
Set-TextValue 'D12' '46.74'
$ws.Range('E12').Value = '  +4.69%  '
Set-TextValue 'D13' '0.0000270'
$ws.Range('E13').Value = '  +0.72%  '
Set-TextValue 'D14' '691.98'
$ws.Range('E14').Value = '  +4.75%  '
$ws.Range('D15').Value = '3.865.35'
$ws.Range('E15').Value = '  +2.39%  '
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').Value = '67.616.22'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').Value = '3.326.94'
Set-TextValue 'D20' '17.46'
$ws.Range('E20').Value = '  +1.58%  '
Set-TextValue 'D21' '11.06'
$ws.Range('E21').Value = '  +3.20%  '
$ws.Range('E22').Value = '  +1.50%  '
Set-TextValue 'D23' '5.45'
$ws.Range('E23').Value = '  +3.06%  '
Set-TextValue 'D24' '16.79'
$ws.Range('E24').Value = '  -0.46%  '
Set-TextValue 'D25' '101.25'
$ws.Range('E25').Value = '  +4.50%  '
$ws.Range('E26').Value = '  +1.59%  '
$ws.Range('E27').Value = '  +1.63%  '
$ws.Range('E28').Value = '  +3.35%  '
Set-TextValue 'D29' '32.81'
$ws.Range('E29').Value = '  +2.66%  '
Set-TextValue 'D30' '8.48'
$ws.Range('E30').Value = '  +2.37%  '
Set-TextValue 'D31' '6.98'
$ws.Range('E31').Value = '  +2.74%  '
Set-TextValue 'D32' '567.51'
$ws.Range('E32').Value = '  +0.22%  '
Set-TextValue 'D33' '10.95'
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('E34').Value = '  +3.03%  '
Set-TextValue 'D35' '0.999'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '3.712.83'
$ws.Range('E36').Value = '  -0.99%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D37' '57.20'
$ws.Range('E37').Value = '  +3.17%  '
Set-TextValue 'D38' '3.24'
$ws.Range('E38').Value = '  -4.83%  '
Set-TextValue 'D39' '34.96'
$ws.Range('E39').Value = '  +8.56%  '
$ws.Range('E40').Value = '  +2.98%  '
$ws.Range('E41').Value = '  +5.09%  '
Set-TextValue 'D42' '2.59'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D43' '0.333'
$ws.Range('E43').Value = '  +2.86%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D44' '3.30'
$ws.Range('E44').Value = '  +2.84%  '
$ws.Range('D45').Value = '0.0₃0664'
$ws.Range('E45').Value = '  +1.42%  '
$ws.Range('E46').Value = '  +2.23%  '
$ws.Range('E47').Value = '  +3.53%  '
$ws.Range('E48').Value = '  +1.58%  '
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('E50').Value = '  -0.85%  '
Set-TextValue 'D51' '131.05'
$ws.Range('E51').Value = '  +2.17%  '
